$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.Goto($ws.Range("A18"), $true)
$ws.Range("C41").Select()
